$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1 ("Sheet1" -> "העברות"): rename only, keep its existing data
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "העברות"

# ------------------------------------------------------------------
# Sheet 2 ("מוסד"): brand new sheet, inserted right after sheet 1,
# holding institution-import data.
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "מוסד"

# Headers (row 1)
$ws2.Range("A1").Value = "מזהה מוסד"
$ws2.Range("B1").Value = "מזהה מוסד שולח"
$ws2.Range("C1").Value = "שם מוסד"
$ws2.Range("D1").NumberFormat = "@"
$ws2.Range("D1").Value = "מספר מזהה"

# Data (row 2) - write D2 before C2 so the shared-string table order
# matches the source data entry order (מזהה -> 001, then שם מוסד).
$ws2.Range("A2").Value = 12345678
$ws2.Range("B2").Value = 12345
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "001"
$ws2.Range("C2").Value = "החברה שלי בע""מ"

# Column widths (autofit-style custom widths)
$ws2.Columns.Item(1).ColumnWidth = 20.944010416666668
$ws2.Columns.Item(2).ColumnWidth = 14.276041666666666
$ws2.Columns.Item(3).ColumnWidth = 19.721354166666668
$ws2.Columns.Item(4).ColumnWidth = 10.830729166666666

$ws2.Range("B5").Select() | Out-Null

# Page setup (matches the default print settings Excel writes for a
# freshly created sheet: Letter/A4-class "paperSize 9", portrait)
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# Back on sheet 1: column widths + selection tweaks
# ------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 9.276041666666666
$ws1.Columns.Item(5).ColumnWidth = 10.385416666666666
$ws1.Range("E10").Select() | Out-Null

# Keep sheet 1 ("העברות") the active / visible tab, matching the source
$ws1.Activate() | Out-Null
